$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 29
$ws_ALC.Range("H29").Value = 1493
$ws_ALC.Range("I29").Value = 1493
$ws_ALC.Range("K29").Value = 4479
$ws_ALC.Range("M29").Value = -4198

# ALC row 129
$ws_ALC.Range("H129").Value = 323616.56
$ws_ALC.Range("I129").Value = 222.625
$ws_ALC.Range("K129").Value = 667.875
$ws_ALC.Range("M129").Value = 4332.125

# ARM row 2
$ws_ARM.Range("H2").Value = 736.84
$ws_ARM.Range("I2").Value = 760.0454999999999
$ws_ARM.Range("J2").Value = 566.6667
$ws_ARM.Range("K2").Value = 760.0454999999999
$ws_ARM.Range("L2").Value = 566.6667
$ws_ARM.Range("M2").Value = -647.0454999999999
$ws_ARM.Range("N2").Value = -792.6667

# ARM row 61
$ws_ARM.Range("H61").Value = 1454.6
$ws_ARM.Range("I61").Value = 1363.4445
$ws_ARM.Range("J61").Value = 2275
$ws_ARM.Range("K61").Value = 1363.4445
$ws_ARM.Range("L61").Value = 2275
$ws_ARM.Range("M61").Value = -1151.4445
$ws_ARM.Range("N61").Value = -2699

# ARM row 116
$ws_ARM.Range("H116").Value = 736.84
$ws_ARM.Range("I116").Value = 760.0454999999999
$ws_ARM.Range("J116").Value = 566.6667
$ws_ARM.Range("K116").Value = 760.0454999999999
$ws_ARM.Range("L116").Value = 566.6667
$ws_ARM.Range("M116").Value = 1533.9545
$ws_ARM.Range("N116").Value = -5154.6667

# ARM row 136
$ws_ARM.Range("H136").Value = 1454.6
$ws_ARM.Range("I136").Value = 1363.4445
$ws_ARM.Range("J136").Value = 2275
$ws_ARM.Range("K136").Value = 4090.3335
$ws_ARM.Range("L136").Value = 6825
$ws_ARM.Range("M136").Value = -1540.3335
$ws_ARM.Range("N136").Value = -11925

# BSM row 3
$ws_BSM.Range("H3").Value = 736.84
$ws_BSM.Range("I3").Value = 760.0454999999999
$ws_BSM.Range("J3").Value = 566.6667
$ws_BSM.Range("K3").Value = 760.0454999999999
$ws_BSM.Range("L3").Value = 566.6667
$ws_BSM.Range("M3").Value = -646.0454999999999
$ws_BSM.Range("N3").Value = -794.6667

# BSM row 20
$ws_BSM.Range("H20").Value = 2310.15
$ws_BSM.Range("I20").Value = 2046.2667
$ws_BSM.Range("J20").Value = 3101.8
$ws_BSM.Range("K20").Value = 2046.2667
$ws_BSM.Range("L20").Value = 3101.8
$ws_BSM.Range("M20").Value = -1799.2667
$ws_BSM.Range("N20").Value = -3595.8

# CRP row 16
$ws_CRP.Range("H16").Value = 1512.3334
$ws_CRP.Range("I16").Value = 1177.75
$ws_CRP.Range("J16").Value = 1780
$ws_CRP.Range("K16").Value = 1177.75
$ws_CRP.Range("L16").Value = 1780
$ws_CRP.Range("M16").Value = -890.75
$ws_CRP.Range("N16").Value = -2354

# CRP row 31
$ws_CRP.Range("H31").Value = 3534.5527
$ws_CRP.Range("I31").Value = 2882.4285
$ws_CRP.Range("J31").Value = 3914.9583
$ws_CRP.Range("K31").Value = 2882.4285
$ws_CRP.Range("L31").Value = 3914.9583
$ws_CRP.Range("M31").Value = -2587.4285
$ws_CRP.Range("N31").Value = -4504.9583

# CRP row 34
$ws_CRP.Range("H34").Value = 3534.5527
$ws_CRP.Range("I34").Value = 2882.4285
$ws_CRP.Range("J34").Value = 3914.9583
$ws_CRP.Range("K34").Value = 2882.4285
$ws_CRP.Range("L34").Value = 3914.9583
$ws_CRP.Range("M34").Value = -2680.4285
$ws_CRP.Range("N34").Value = -4318.9583

# CRP row 99
$ws_CRP.Range("H99").Value = 3509.95
$ws_CRP.Range("I99").Value = 2857.1428
$ws_CRP.Range("K99").Value = 2857.1428
$ws_CRP.Range("M99").Value = -1359.1428

# CRP row 113
$ws_CRP.Range("H113").Value = 1512.3334
$ws_CRP.Range("I113").Value = 1177.75
$ws_CRP.Range("J113").Value = 1780
$ws_CRP.Range("K113").Value = 1177.75
$ws_CRP.Range("L113").Value = 1780
$ws_CRP.Range("M113").Value = 992.25
$ws_CRP.Range("N113").Value = -6120

# CRP row 126
$ws_CRP.Range("H126").Value = 3509.95
$ws_CRP.Range("I126").Value = 2857.1428
$ws_CRP.Range("K126").Value = 8571.428400000001
$ws_CRP.Range("M126").Value = -6101.428400000001

# CRP row 131
$ws_CRP.Range("H131").Value = 0
$ws_CRP.Range("J131").Value = 0
$ws_CRP.Range("L131").Value = 0
$ws_CRP.Range("N131").ClearContents()

# CUL row 5
$ws_CUL.Range("H5").Value = 1951.5555
$ws_CUL.Range("I5").Value = 1163.5
$ws_CUL.Range("J5").Value = 2582
$ws_CUL.Range("K5").Value = 3490.5
$ws_CUL.Range("L5").Value = 7746
$ws_CUL.Range("M5").Value = -3378.5
$ws_CUL.Range("N5").Value = -7970

# CUL row 113
$ws_CUL.Range("H113").Value = 637.3158
$ws_CUL.Range("I113").Value = 512.8
$ws_CUL.Range("J113").Value = 775.6667
$ws_CUL.Range("K113").Value = 1538.4
$ws_CUL.Range("L113").Value = 2327.0001
$ws_CUL.Range("M113").Value = 631.6000000000001
$ws_CUL.Range("N113").Value = -6667.0001

# CUL row 117
$ws_CUL.Range("H117").Value = 1403.625
$ws_CUL.Range("I117").Value = 1257.25
$ws_CUL.Range("J117").Value = 1550
$ws_CUL.Range("K117").Value = 3771.75
$ws_CUL.Range("L117").Value = 4650
$ws_CUL.Range("M117").Value = -329.75
$ws_CUL.Range("N117").Value = -11534

# CUL row 131
$ws_CUL.Range("H131").Value = 725.23
$ws_CUL.Range("J131").Value = 729.2222
$ws_CUL.Range("L131").Value = 2187.6666
$ws_CUL.Range("N131").Value = -12267.6666

# CUL row 135
$ws_CUL.Range("H135").Value = 1951.5555
$ws_CUL.Range("I135").Value = 1163.5
$ws_CUL.Range("J135").Value = 2582
$ws_CUL.Range("K135").Value = 10471.5
$ws_CUL.Range("L135").Value = 23238
$ws_CUL.Range("M135").Value = -7936.5
$ws_CUL.Range("N135").Value = -28308

# GSM row 123
$ws_GSM.Range("H123").Value = 5839.3
$ws_GSM.Range("I123").Value = 3228.75
$ws_GSM.Range("J123").Value = 16281.5
$ws_GSM.Range("K123").Value = 3228.75
$ws_GSM.Range("L123").Value = 16281.5
$ws_GSM.Range("M123").Value = -778.75
$ws_GSM.Range("N123").Value = -21181.5

# LTW row 7
$ws_LTW.Range("H7").Value = 2251.0386
$ws_LTW.Range("I7").Value = 2141.15
$ws_LTW.Range("J7").Value = 2617.3333
$ws_LTW.Range("K7").Value = 2141.15
$ws_LTW.Range("L7").Value = 2617.3333
$ws_LTW.Range("M7").Value = -2029.15
$ws_LTW.Range("N7").Value = -2841.3333

# LTW row 88
$ws_LTW.Range("H88").Value = 0
$ws_LTW.Range("J88").Value = 0
$ws_LTW.Range("L88").Value = 0
$ws_LTW.Range("N88").ClearContents()

# LTW row 91
$ws_LTW.Range("H91").Value = 0
$ws_LTW.Range("J91").Value = 0
$ws_LTW.Range("L91").Value = 0
$ws_LTW.Range("N91").ClearContents()

# LTW row 126
$ws_LTW.Range("H126").Value = 2251.0386
$ws_LTW.Range("I126").Value = 2141.15
$ws_LTW.Range("J126").Value = 2617.3333
$ws_LTW.Range("K126").Value = 6423.450000000001
$ws_LTW.Range("L126").Value = 7851.999899999999
$ws_LTW.Range("M126").Value = -3953.450000000001
$ws_LTW.Range("N126").Value = -12791.9999

# LTW row 136
$ws_LTW.Range("H136").Value = 1464.9286
$ws_LTW.Range("I136").Value = 1467.5
$ws_LTW.Range("J136").Value = 1449.5
$ws_LTW.Range("K136").Value = 4402.5
$ws_LTW.Range("L136").Value = 4348.5
$ws_LTW.Range("M136").Value = -1852.5
$ws_LTW.Range("N136").Value = -9448.5

# WVR row 54
$ws_WVR.Range("H54").Value = 13333.333
$ws_WVR.Range("J54").Value = 13333.333
$ws_WVR.Range("L54").Value = 13333.333
$ws_WVR.Range("N54").Value = -14373.333

# WVR row 126
$ws_WVR.Range("H126").Value = 1118.4615
$ws_WVR.Range("I126").Value = 1190.4783
$ws_WVR.Range("J126").Value = 566.3333
$ws_WVR.Range("K126").Value = 3571.4349
$ws_WVR.Range("L126").Value = 1698.9999
$ws_WVR.Range("M126").Value = -1101.4349
$ws_WVR.Range("N126").Value = -6638.9999
